# Update usmap_transform() to account for Puerto Rico.
#
# The "citypop" sheet lists one row per state (lon, lat, state, abbr,
# most_populous_city, city_pop), sorted by state. Puerto Rico needs to be
# inserted as its own row, in order, between Pennsylvania (row 40) and
# Rhode Island (old row 41) - i.e. it becomes the new row 41, and every
# row from the old row 41 onward shifts down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 41; this pushes the old rows 41-52 down to 42-53.
$ws.Rows.Item(41).Insert()

# The newly inserted row has no formatting yet. Copy the look (number
# formats / fonts / borders / row height) from the row directly below it
# (old row 41, now at 42) so the new row matches the rest of the table.
$ws.Range("A42:F42").Copy()
$ws.Range("A41:F41").PasteSpecial(-4122)
$ws.Rows.Item(41).RowHeight = $ws.Rows.Item(42).RowHeight

# Fill in the Puerto Rico data.
$ws.Range("A41").Value = -66.06
$ws.Range("B41").Value = 18.41
$ws.Range("C41").Value = "Puerto Rico"
$ws.Range("D41").Value = "PR"
$ws.Range("E41").Value = "San Juan"
$ws.Range("F41").Value = 395326

Write-Host "Inserted Puerto Rico as row 41 (citypop sheet now has 53 data rows)"
